# Generate Report for Handback
#
# The localization handback for af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md has
# completed (it's now "Handed back: in sync with en-US" instead of still
# "Ready for handoff"), so the localization-status report is refreshed to
# reflect that for both target languages (zh-cn, de-de) and the Overview
# rollup sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: the af3f0b90 row (row 3) now shows both locales as
#     handed back, in sync with en-US. ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: af3f0b90 row (row 3) - status flips to handed back, the
#     handback timestamp is recorded, and the stale error detail clears. ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-30 19:01:23"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: af3f0b90 row (row 3) - same treatment, different
#     handback timestamp. ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-30 19:01:31"
$wsDeDe.Range("P3").Value = ""
